$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the two new columns (I0, IF), matching the header style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Fill data rows 2-31: I = 1 (constant), J = copy of H value
for ($r = 2; $r -le 31; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
